# Prefix each worksheet's own name (plus a space) onto the existing
# Step/command names found in column A (rows 2..last used row) for the
# "protocol" sheets that hold Name/Text/Note/*Guidelines tables.
#
# This mirrors the commit: "fix: unique command names in XLSX - prefix
# protocol name to each step"

$wb = $excel.ActiveWorkbook

# Sheets whose column-A labels need the sheet name prefixed.
$targetSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the last used row based on column A.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($current -ne $null -and $current -ne "") {
            $prefix = "$sheetName "
            if (-not $current.ToString().StartsWith($prefix)) {
                $cell.Value = "$sheetName $current"
            }
        }
    }
}
